# Add SVR parameter loading from pred_par structure and Excel files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) for SVR parameters
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New data cells (row 2) with the SVR parameter values
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# Move the active selection to J9, matching the author's final cursor position
$ws.Range("J9").Select()
